# Add a 4th column (D) with "Time(Just Search, 4 Thread) in Release" timings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column.
$ws.Range("D6").Value = "Time(Just Search, 4 Thread) in Release"

# Average formula for the new column (same pattern as column C).
$ws.Range("D7").Formula = "=AVERAGE(D8:D107)"

# 100 data points for rows 8-107.
$dValues = @(0.38400000000000001,0.311,0.28899999999999998,0.27400000000000002,0.28100000000000003,0.28100000000000003,0.26600000000000001,0.27800000000000002,0.28000000000000003,0.29199999999999998,0.26400000000000001,0.25700000000000001,0.26600000000000001,0.26800000000000002,0.247,0.248,0.27500000000000002,0.26400000000000001,0.25700000000000001,0.25800000000000001,0.26200000000000001,0.251,0.26200000000000001,0.25700000000000001,0.25900000000000001,0.26,0.252,0.26300000000000001,0.26800000000000002,0.26200000000000001,0.252,0.254,0.25600000000000001,0.26100000000000001,0.312,0.26300000000000001,0.26600000000000001,0.27900000000000003,0.26300000000000001,0.26500000000000001,0.28100000000000003,0.30099999999999999,0.28499999999999998,0.29699999999999999,0.28499999999999998,0.28100000000000003,0.28999999999999998,0.28799999999999998,0.28499999999999998,0.29599999999999999,0.27700000000000002,0.27100000000000002,0.28599999999999998,0.28699999999999998,0.28199999999999997,0.27300000000000002,0.28299999999999997,0.29499999999999998,0.34799999999999998,0.29599999999999999,0.27300000000000002,0.28999999999999998,0.29399999999999998,0.29099999999999998,0.28599999999999998,0.317,0.28399999999999997,0.28799999999999998,0.28599999999999998,0.27600000000000002,0.27600000000000002,0.29199999999999998,0.28499999999999998,0.28799999999999998,0.28799999999999998,0.28899999999999998,0.29399999999999998,0.29899999999999999,0.28199999999999997,0.28199999999999997,0.28399999999999997,0.29799999999999999,0.27400000000000002,0.27300000000000002,0.27800000000000002,0.28599999999999998,0.28599999999999998,0.29199999999999998,0.28699999999999998,0.28299999999999997,0.28599999999999998,0.28899999999999998,0.28999999999999998,0.27800000000000002,0.28499999999999998,0.29799999999999999,0.29099999999999998,0.28899999999999998,0.27500000000000002,0.28199999999999997)

for ($i = 0; $i -lt $dValues.Length; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}

# Match the column D width used in the edited workbook as closely as this
# engine's pixel-snapped ColumnWidth model allows.
$ws.Columns.Item(4).ColumnWidth = 34.5

# Move the active selection to D8, matching the post-edit selection.
$ws.Range("D8").Select()
